$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

# New user rows 27-30: column A gets a new unique code, column C reuses the
# "001" value/format already used e.g. in row 2 / row 26 (style index 4,
# text-formatted number).
$newCodes = @("F02062", "F03153", "F00191", "F02729")

for ($i = 0; $i -lt $newCodes.Length; $i++) {
    $row = 27 + $i

    # Column A: plain text value (shared string, no explicit style).
    $ws.Cells.Item($row, 1).Value = $newCodes[$i]

    # Column C: copy formatting from the row above (keeps it on style index 4)
    # then set the text value "001" (reuses existing shared string).
    $ws.Cells.Item($row - 1, 3).Copy() | Out-Null
    $ws.Cells.Item($row, 3).PasteSpecial(-4122) | Out-Null
    $ws.Cells.Item($row, 3).Value = "001"
}

$excel.CutCopyMode = 0

# Update selection to match the post-edit state.
$ws.Range("B26").Select() | Out-Null
